$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header for column C
$ws.Range("C1").Value = "timestamp"

# Update row 2 values
$ws.Range("A2").Value = "老師A"
$ws.Range("B2").Value = "今天辛苦了！"
$ws.Range("C2").Value = "2025-07-22 15:00"

# Remove row 3 entirely (shrinks used range to A1:C2)
$ws.Rows.Item(3).Delete()
